$d = $word.ActiveDocument
$full = $d.Range().WordOpenXML
$full = $full -replace 'w14:paraId="[0-9A-Fa-f]+" w14:textId="[0-9A-Fa-f]+" ', ''

# --- Block 1: replace paragraphs from "Celność" marker through end of "Ofensywny" paragraph ---
$marker1Start = '<w:p w:rsidR="00494EBF"'
$marker1End = '<w:p w:rsidR="004A18BE"'
$s1 = $full.IndexOf($marker1Start)
$e1 = $full.IndexOf($marker1End)
if ($s1 -lt 0 -or $e1 -lt 0) { Write-Host "ERROR: block1 markers not found: $s1 $e1" }
$new1 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Pancerz:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Naturalny i nabyty?</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Czy zadawanie obrażeń = zmiana w pancerzu i odpornościach</w:t></w:r><w:r><w:t xml:space="preserve"> (woda, ziemia, powietrze, </w:t></w:r><w:r><w:t>ogień</w:t></w:r><w:r><w:t xml:space="preserve"> + pominięcie odporności z danej rasy</w:t></w:r><w:r><w:t>)</w:t></w:r><w:r><w:t>?</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Nisczenie pancerza? Nieodwracalne?</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Celność  - co z nią zrobić?</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>3</w:t></w:r><w:r><w:t xml:space="preserve"> tryby akcji na polu bitwy:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Pass</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Ofensywny</w:t></w:r><w:r><w:t xml:space="preserve"> – szansa na zadanie obrażeń krytycznych,</w:t></w:r></w:p>'
$full = $full.Substring(0, $s1) + $new1 + $full.Substring($e1)

# --- Block 2: replace paragraphs from "Do zrobienia" marker through end of following empty paragraph ---
$marker2Start = '<w:p w:rsidR="00337B37"'
$marker2End = '<w:sectPr'
$s2 = $full.IndexOf($marker2Start)
$e2 = $full.IndexOf($marker2End)
if ($s2 -lt 0 -or $e2 -lt 0) { Write-Host "ERROR: block2 markers not found: $s2 $e2" }
$new2 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Do zrobienia</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Model </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t>Akcji:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>AttributeChange</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>[] Buff</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Pętla - </w:t></w:r><w:r><w:t>Mechanizm atakowania</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Zapisz id obecnie wybranej jednostki</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Wybierz jednostkę do zaatakowania</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Zapisz jej id</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Wyślij atakującemu #ref do jednostki</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Wyczyść stan gry</w:t></w:r></w:p>'
$full = $full.Substring(0, $s2) + $new2 + $full.Substring($e2)

# --- numbering.xml: remove tentative="1" on the w:ilvl="3" w:tplc="04150001" level ---
$old3 = '<w:lvl w:ilvl="3" w:tplc="04150001" w:tentative="1">'
$new3 = '<w:lvl w:ilvl="3" w:tplc="04150001">'
if ($full.IndexOf($old3) -lt 0) { Write-Host "ERROR: old3 not found" }
$full = $full.Replace($old3, $new3)

$d.Range().InsertXML($full)
Write-Host "DONE"
